$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf2"
$ws.Range("C2").Value = "Nrp1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.6462393333333333
$ws.Range("H2").Value = 1.938718
$ws.Range("I2").Value = 0.03461850536298827
$ws.Range("J2").Value = 0.03461850536298827
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 102.8289443333334
$ws.Range("N2").Value = 308.486833
$ws.Range("O2").Value = 0.5559120396302444
$ws.Range("P2").Value = 0.5559120396302443
$ws.Range("Q2").Value = 66.45210843334378
$ws.Range("R2").Value = 598.0689759000941
$ws.Range("S2").Value = 0.01924484392528936
$ws.Range("T2").Value = 0.01924484392528936

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf2"
$ws.Range("C3").Value = "Nrp1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.6462393333333333
$ws.Range("H3").Value = 1.938718
$ws.Range("I3").Value = 0.03461850536298827
$ws.Range("J3").Value = 0.03461850536298827
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 63.66262833333334
$ws.Range("N3").Value = 190.987885
$ws.Range("O3").Value = 0.3441717873742006
$ws.Range("P3").Value = 0.3441717873742006
$ws.Range("Q3").Value = 41.14129449238111
$ws.Range("R3").Value = 370.27165043143
$ws.Range("S3").Value = 0.01191471286700302
$ws.Range("T3").Value = 0.01191471286700302

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf2"
$ws.Range("C4").Value = "Nrp1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.6462393333333333
$ws.Range("H4").Value = 1.938718
$ws.Range("I4").Value = 0.03461850536298827
$ws.Range("J4").Value = 0.03461850536298827
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 18.481835
$ws.Range("N4").Value = 55.445505
$ws.Range("O4").Value = 0.09991617299555507
$ws.Range("P4").Value = 0.09991617299555505
$ws.Range("Q4").Value = 11.94368872917667
$ws.Range("R4").Value = 107.49319856259
$ws.Range("S4").Value = 0.003458948570695887
$ws.Range("T4").Value = 0.003458948570695886

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf2"
$ws.Range("C5").Value = "Nrp1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 15.322826
$ws.Range("H5").Value = 45.968478
$ws.Range("I5").Value = 0.8208310864042159
$ws.Range("J5").Value = 0.8208310864042158
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 102.8289443333334
$ws.Range("N5").Value = 308.486833
$ws.Range("O5").Value = 0.5559120396302444
$ws.Range("P5").Value = 0.5559120396302443
$ws.Range("Q5").Value = 1575.630021783353
$ws.Range("R5").Value = 14180.67019605018
$ws.Range("S5").Value = 0.456309883434877
$ws.Range("T5").Value = 0.4563098834348768

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf2"
$ws.Range("C6").Value = "Nrp1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 15.322826
$ws.Range("H6").Value = 45.968478
$ws.Range("I6").Value = 0.8208310864042159
$ws.Range("J6").Value = 0.8208310864042158
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 63.66262833333334
$ws.Range("N6").Value = 190.987885
$ws.Range("O6").Value = 0.3441717873742006
$ws.Range("P6").Value = 0.3441717873742006
$ws.Range("Q6").Value = 975.4913766543368
$ws.Range("R6").Value = 8779.42238988903
$ws.Range("S6").Value = 0.2825069021400459
$ws.Range("T6").Value = 0.2825069021400459

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf2"
$ws.Range("C7").Value = "Nrp1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 15.322826
$ws.Range("H7").Value = 45.968478
$ws.Range("I7").Value = 0.8208310864042159
$ws.Range("J7").Value = 0.8208310864042158
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 18.481835
$ws.Range("N7").Value = 55.445505
$ws.Range("O7").Value = 0.09991617299555507
$ws.Range("P7").Value = 0.09991617299555505
$ws.Range("Q7").Value = 283.19394186571
$ws.Range("R7").Value = 2548.74547679139
$ws.Range("S7").Value = 0.08201430082929304
$ws.Range("T7").Value = 0.08201430082929302

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fgf2"
$ws.Range("C8").Value = "Nrp1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.698388
$ws.Range("H8").Value = 8.095164
$ws.Range("I8").Value = 0.1445504082327959
$ws.Range("J8").Value = 0.1445504082327959
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 102.8289443333334
$ws.Range("N8").Value = 308.486833
$ws.Range("O8").Value = 0.5559120396302444
$ws.Range("P8").Value = 0.5559120396302443
$ws.Range("Q8").Value = 277.4723894417347
$ws.Range("R8").Value = 2497.251504975613
$ws.Range("S8").Value = 0.08035731227007803
$ws.Range("T8").Value = 0.08035731227007802

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fgf2"
$ws.Range("C9").Value = "Nrp1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.698388
$ws.Range("H9").Value = 8.095164
$ws.Range("I9").Value = 0.1445504082327959
$ws.Range("J9").Value = 0.1445504082327959
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 63.66262833333334
$ws.Range("N9").Value = 190.987885
$ws.Range("O9").Value = 0.3441717873742006
$ws.Range("P9").Value = 0.3441717873742006
$ws.Range("Q9").Value = 171.7864723431267
$ws.Range("R9").Value = 1546.07825108814
$ws.Range("S9").Value = 0.04975017236715172
$ws.Range("T9").Value = 0.04975017236715172

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fgf2"
$ws.Range("C10").Value = "Nrp1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.698388
$ws.Range("H10").Value = 8.095164
$ws.Range("I10").Value = 0.1445504082327959
$ws.Range("J10").Value = 0.1445504082327959
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 18.481835
$ws.Range("N10").Value = 55.445505
$ws.Range("O10").Value = 0.09991617299555507
$ws.Range("P10").Value = 0.09991617299555505
$ws.Range("Q10").Value = 49.87116178198
$ws.Range("R10").Value = 448.84045603782
$ws.Range("S10").Value = 0.01444292359556614
$ws.Range("T10").Value = 0.01444292359556614
